$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Fecha(D), Variedad(H), Calidad(I), Volumen(J),
#             Precio minimo(K), Precio maximo(L), Precio promedio ponderado(M),
#             Origen(O), Precio $/Kg(P)
$data = @(
  @(363, 45021, "Calameño", "Primera", 60,  1500, 1500, 1500, "Región del Maule",    1500),
  @(364, 45021, "Calameño", "Segunda", 40,  1000, 1000, 1000, "Región del Maule",    1000),
  @(365, 45021, "Tuna",     "Primera", 80,  1300, 1300, 1300, "Región del Maule",    1300),
  @(366, 45021, "Tuna",     "Segunda", 90,  800,  1000, 933,  "Región del Maule",    933),
  @(367, 45002, "Plátano",  "Segunda", 500, 1500, 1500, 1500, "Región del Maule",    1500),
  @(368, 44970, "Calameño", "Extra",   500, 900,  900,  900,  "Región del Maule",    900),
  @(369, 44970, "Calameño", "Primera", 500, 700,  700,  700,  "Región del Maule",    700),
  @(370, 44970, "Calameño", "Segunda", 500, 500,  500,  500,  "Región del Maule",    500),
  @(371, 44970, "Tuna",     "Extra",   500, 900,  900,  900,  "Región del Maule",    900),
  @(372, 44970, "Tuna",     "Primera", 500, 700,  700,  700,  "Región del Maule",    700),
  @(373, 44970, "Tuna",     "Segunda", 500, 500,  500,  500,  "Región del Maule",    500),
  @(374, 44172, "Calameño", "Segunda", 300, 900,  1000, 950,  "Región de O'Higgins", 950),
  @(375, 44172, "Calameño", "Tercera", 160, 800,  850,  825,  "Región de O'Higgins", 825),
  @(376, 44952, "Calameño", "Extra",   500, 1000, 1000, 1000, "Región del Maule",    1000),
  @(377, 44952, "Calameño", "Primera", 500, 800,  800,  800,  "Región del Maule",    800),
  @(378, 44952, "Calameño", "Segunda", 500, 600,  600,  600,  "Región del Maule",    600),
  @(379, 44952, "Tuna",     "Extra",   500, 1000, 1000, 1000, "Región del Maule",    1000),
  @(380, 44952, "Tuna",     "Segunda", 500, 600,  600,  600,  "Región del Maule",    600),
  @(381, 44252, "Tuna",     "Primera", 400, 900,  950,  925,  "Región del Maule",    925),
  @(382, 44252, "Tuna",     "Segunda", 600, 800,  850,  825,  "Región del Maule",    825),
  @(383, 44243, "Calameño", "Primera", 600, 700,  750,  725,  "Región del Maule",    725),
  @(384, 44243, "Calameño", "Segunda", 500, 600,  650,  625,  "Región del Maule",    625),
  @(385, 45007, "Calameño", "Primera", 200, 1500, 1500, 1500, "Región del Maule",    1500),
  @(386, 45007, "Calameño", "Segunda", 100, 1000, 1000, 1000, "Región del Maule",    1000),
  @(387, 45007, "Calameño", "Tercera", 80,  500,  500,  500,  "Región del Maule",    500)
)

foreach ($row in $data) {
    $r = $row[0]

    # Columns that are constant across the whole block (A..C, E..G, N, Q, R).
    # Only strictly necessary for the brand-new rows (384-387), but harmless
    # to (re)write for every row.
    $ws.Cells.Item($r, 1).Value = 7
    $ws.Cells.Item($r, 2).Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Cells.Item($r, 3).Value = "Ñuble"

    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 4).Value = $row[1]

    $ws.Cells.Item($r, 5).Value = 16
    $ws.Cells.Item($r, 6).Value = 100112027
    $ws.Cells.Item($r, 7).Value = "Melón"

    $ws.Cells.Item($r, 8).Value = $row[2]
    $ws.Cells.Item($r, 9).Value = $row[3]
    $ws.Cells.Item($r, 10).Value = $row[4]
    $ws.Cells.Item($r, 11).Value = $row[5]
    $ws.Cells.Item($r, 12).Value = $row[6]
    $ws.Cells.Item($r, 13).Value = $row[7]

    $ws.Cells.Item($r, 14).Value = "$/unidad"
    $ws.Cells.Item($r, 15).Value = $row[8]
    $ws.Cells.Item($r, 16).Value = $row[9]
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
